$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.568.87"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +0.52%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.819.67"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +1.23%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.009"
$cell.ClearFormats()
$ws.Range("E4").Value = "  +0.15%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "305.42"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -0.66%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4663"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +2.07%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3589"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -0.89%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07121"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +0.08%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.8980"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +1.65%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07778"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -0.58%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "19.30"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -1.22%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.834.57"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +6.80%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.241"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -0.89%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "6.325"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -0.14%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "87.26"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +2.46%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "1.010"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +0.16%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000008538"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("E19").Value = "  +0.20%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "26.613.11"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +0.58%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.14"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -1.01%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.003"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("E23").Value = "  +0.08%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.916"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -3.22%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "151.93"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -0.42%  "

$ws.Range("E26").Value = "  -0.32%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.970"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -3.86%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "113.44"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.23%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "4.798"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -1.57%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.08793"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +1.51%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.131"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +2.53%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.7279"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.720"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.426"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -0.67%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.120"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("E36").Value = "  -0.25%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.01920"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -1.39%  "

$ws.Range("E38").Value = "  +1.50%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.05092"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -0.73%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "6.822"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -1.18%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.5032"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -3.13%  "

$ws.Range("E42").Value = "  -2.70%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "7.940"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -1.16%  "

$ws.Range("E44").Value = "  +0.22%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.4645"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.86%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "9.913"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -0.17%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "97.50"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -2.97%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.553"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -2.37%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.05980"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +0.08%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "63.52"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -1.27%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "35.80"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -1.33%  "
